# Regenerate the "K" column (column G) values for the save_data sheet.
# Per the commit message: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" — the recalculated K values for
# each existing data row (rows 2-48) are written below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 0
    35 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 0
    43 = 1
    45 = 2
    46 = 0
    47 = 1
    48 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
